# This script applies the change described by the diff:
# - The "catch" row (G, J, K, L only) that used to live at row 16
#   (stimuli/catch_17_stairs.jpg) is replaced by a new catch row
#   at row 12 with stimulus "stimuli/catch_20.jpg".
# - The four "target" rows that used to occupy rows 12-15
#   (img_bpyv5.png, img_ozxpp.png, img_wijef.png, img_a8y4y.png)
#   each shift down by one row, to rows 13-16 respectively.
# Columns A-F (subject_id, task, block_total, block_scene, trial_block,
# trial_total) are left untouched since they already increment per-row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: capture current (pre-edit) values for rows 12-15 (H..V) ---
# These values will be "pushed down" into rows 13-16.
$rowsToShift = @(12, 13, 14, 15)
$captured = @{}

foreach ($r in $rowsToShift) {
    $captured[$r] = @{
        H = $ws.Range("H$r").Value2
        I = $ws.Range("I$r").Value2
        J = $ws.Range("J$r").Value2
        K = $ws.Range("K$r").Value2
        L = $ws.Range("L$r").Value2
        M = $ws.Range("M$r").Value2
        N = $ws.Range("N$r").Value2
        O = $ws.Range("O$r").Value2
        P = $ws.Range("P$r").Value2
        Q = $ws.Range("Q$r").Value2
        R = $ws.Range("R$r").Value2
        S = $ws.Range("S$r").Value2
        T = $ws.Range("T$r").Value2
        U = $ws.Range("U$r").Value2
        V = $ws.Range("V$r").Value2
    }
}

# --- Step 2: write the captured values down one row (12->13, 13->14, 14->15, 15->16) ---
foreach ($r in $rowsToShift) {
    $target = $r + 1
    $vals = $captured[$r]

    $ws.Range("H$target").Value = $vals.H
    $ws.Range("I$target").Value = $vals.I
    $ws.Range("J$target").Value = $vals.J
    $ws.Range("K$target").Value = $vals.K
    $ws.Range("L$target").Value = $vals.L
    $ws.Range("M$target").Value = $vals.M
    $ws.Range("N$target").Value = $vals.N
    $ws.Range("O$target").Value = $vals.O
    $ws.Range("P$target").Value = $vals.P
    $ws.Range("Q$target").Value = $vals.Q
    $ws.Range("R$target").Value = $vals.R
    $ws.Range("S$target").Value = $vals.S
    $ws.Range("T$target").Value = $vals.T
    $ws.Range("U$target").Value = $vals.U
    $ws.Range("V$target").Value = $vals.V
}

# --- Step 3: turn row 12 into the new catch row ---
# Clear the (now stale/duplicated) target-row-only columns H, I first.
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()

# Set the catch-row fields.
$ws.Range("J12").Value = "catch"
$ws.Range("K12").Value = "f"
$ws.Range("L12").Value = "stimuli/catch_20.jpg"

# Clear the numeric columns M..V on row 12, since a catch row has no
# conceptual/perceptual/typicality/n/p_*/r_* data.
$ws.Range("M12:V12").ClearContents()
